# Insert a new daily record at row 17 of the Papaya / Vega Modelo de Temuco
# sheet. This pushes the existing rows 17-70 down to 18-71 (Excel carries
# the formatting, e.g. the date style on column D, along with the shift),
# growing the sheet's used range from A1:T70 to A1:T71. We then populate
# the freshly inserted row 17 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17:70 down by one row, inserting a blank row at 17.
$ws.Rows("17:17").Insert()

# Populate the new row 17 with the new record.
$ws.Range("A17").Value2 = 10
$ws.Range("B17").Value2 = "Vega Modelo de Temuco"
$ws.Range("C17").Value2 = "La Araucanía"
$ws.Range("D17").Value2 = 44623
$ws.Range("E17").Value2 = 9
$ws.Range("F17").Value2 = "Fruta"
$ws.Range("G17").Value2 = 100108
$ws.Range("H17").Value2 = "Tropicales y subtropicales"
$ws.Range("I17").Value2 = 100108004
$ws.Range("J17").Value2 = "Papaya"
$ws.Range("K17").Value2 = "Cultivar IV Región"
$ws.Range("L17").Value2 = "Primera"
$ws.Range("M17").Value2 = 95
$ws.Range("N17").Value2 = 23000
$ws.Range("O17").Value2 = 23000
$ws.Range("P17").Value2 = 23000
$ws.Range("Q17").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R17").Value2 = "Provincia del Elquí"
$ws.Range("S17").Value2 = 2300
$ws.Range("T17").Value2 = 10
